$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2 through 115
# from 45190 (2023-09-21) to 45192 (2023-09-23).
for ($r = 2; $r -le 115; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}
